$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "COMMENTS" worksheet right after "LECTURES" (so the sheet
#    order becomes STUDENTS, LECTURES, COMMENTS, CHAIRS, ATTENDS).
# ---------------------------------------------------------------------------
$lectures = $wb.Worksheets.Item("LECTURES")
$comments = $wb.Worksheets.Add($null, $lectures)
$comments.Name = "COMMENTS"

# ---------------------------------------------------------------------------
# 2. Headers
# ---------------------------------------------------------------------------
$comments.Range("A1").Value = "ID"
$comments.Range("B1").Value = "STUDENT_ID"
$comments.Range("C1").Value = "LECTURE_ID"
$comments.Range("D1").Value = "TIME"
$comments.Range("E1").Value = "DATE"
$comments.Range("F1").Value = "CONTENT"

# ---------------------------------------------------------------------------
# 3. Data rows (STUDENT_ID, LECTURE_ID, TIME fraction, DATE serial, CONTENT)
# ---------------------------------------------------------------------------
$comments.Range("B2").Value = 3953
$comments.Range("C2").Value = 11
$comments.Range("D2").Value = 0.52222222222222225
$comments.Range("E2").Value = 43256
$comments.Range("F2").Value = "4 out of 7, would recommend"

$comments.Range("B3").Value = 9764
$comments.Range("C3").Value = 11
$comments.Range("D3").Value = 0.5229166666666667
$comments.Range("E3").Value = 43256
$comments.Range("F3").Value = "Great"

$comments.Range("B4").Value = 8064
$comments.Range("C4").Value = 11
$comments.Range("D4").Value = 0.67361111111111116
$comments.Range("E4").Value = 43256
$comments.Range("F4").Value = "Could be better organized"

$comments.Range("B5").Value = 7355
$comments.Range("C5").Value = 11
$comments.Range("D5").Value = 0.71527777777777779
$comments.Range("E5").Value = 43256
$comments.Range("F5").Value = "It's bollocks"

$comments.Range("B6").Value = 971
$comments.Range("C6").Value = 11
$comments.Range("D6").Value = 0.95833333333333337
$comments.Range("E6").Value = 43256
$comments.Range("F6").Value = "Nah"

# ---------------------------------------------------------------------------
# 4. Number formats.
#    Apply each format once, then copy/paste-special(formats) onto the rest
#    of the column so every cell in the column shares a single style record
#    (mirrors how the source workbook only carries 3 new cellXfs entries).
# ---------------------------------------------------------------------------
$comments.Range("B2").NumberFormat = "0"
$comments.Range("B2").Copy()
$comments.Range("B3:B6").PasteSpecial(-4122)

$comments.Range("E2").NumberFormat = "mm-dd-yy"
$comments.Range("E2").Copy()
$comments.Range("E3:E6").PasteSpecial(-4122)

$comments.Range("F2").NumberFormat = "d-mmm"

$comments.Range("D2").NumberFormat = "h:mm"
$comments.Range("D2").Copy()
$comments.Range("D3:D6").PasteSpecial(-4122)

$comments.Application.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. Column widths for the new sheet.
# ---------------------------------------------------------------------------
$comments.Columns.Item(2).ColumnWidth = 11
$comments.Columns.Item(3).ColumnWidth = 11
$comments.Columns.Item(4).ColumnWidth = 11.1
$comments.Columns.Item(5).ColumnWidth = 9.7
$comments.Columns.Item(6).ColumnWidth = 35.9

# ---------------------------------------------------------------------------
# 6. Selections per sheet (restores each sheet's saved cursor position).
# ---------------------------------------------------------------------------
$students = $wb.Worksheets.Item("STUDENTS")
$null = $students.Range("A2:A6").Select()

$null = $lectures.Range("H16").Select()

$chairs = $wb.Worksheets.Item("CHAIRS")
$null = $chairs.Range("E5").Select()

$attends = $wb.Worksheets.Item("ATTENDS")
$null = $attends.Range("D9").Select()

$null = $comments.Range("D7").Select()

# COMMENTS is the sheet that ends up active/selected.
$null = $comments.Activate()
